# Updating BB pathway export: drop the "node" and "node__commodity" sheets,
# rework "demand" (add a Base-scenario column) and "p_unit" (rename/restructure
# the efficiency/operating-point columns), and add the new "unit", "commodity"
# and "p_commodity_price" sheets.

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Drop sheets that are no longer exported.
# ---------------------------------------------------------------------------
[void]$wb.Worksheets.Item("node").Delete()
[void]$wb.Worksheets.Item("node__commodity").Delete()

# ---------------------------------------------------------------------------
# 2. "demand": insert a scenario row under the header and shift the hourly
#    data down by one row.
# ---------------------------------------------------------------------------
$demand = $wb.Worksheets.Item("demand")
$demand.Rows.Item(2).Insert()
$demand.Range("B2").Value = "Base"

# ---------------------------------------------------------------------------
# 3. "p_unit": rename the parameter columns and add the new op00/op01 columns.
# ---------------------------------------------------------------------------
$p_unit = $wb.Worksheets.Item("p_unit")
$p_unit.Range("C1").Value = "eff00"
$p_unit.Range("D1").Value = "eff01"
$p_unit.Range("E1").Value = "op00"
$p_unit.Range("F1").Value = "op01"

$p_unit.Range("C2").Value = 0.4
$p_unit.Range("D2").Value = 0.5
$p_unit.Range("F2").Value = 1

$p_unit.Range("C3").ClearContents()
$p_unit.Range("D3").Value = 0.6

# ---------------------------------------------------------------------------
# 4. Add the new sheets at the end, in order: node, unit, commodity,
#    p_commodity_price.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$node = $wb.Worksheets.Add($null, $lastSheet)
$node.Name = "node"
$node.Range("A1").Value = "gas"
$node.Range("A2").Value = "nodeElec"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$unit = $wb.Worksheets.Add($null, $lastSheet)
$unit.Name = "unit"
$unit.Range("A1").Value = "gas_turbine"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$commodity = $wb.Worksheets.Add($null, $lastSheet)
$commodity.Name = "commodity"
$commodity.Range("A1").Value = "gas"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$priceSheet = $wb.Worksheets.Add($null, $lastSheet)
$priceSheet.Name = "p_commodity_price"
$priceSheet.Range("C1").Value = "gas"

$priceSheet.Range("A2").Value = "2020-01-01T00:00:00"
$priceSheet.Range("B2").Value = "Base"
$priceSheet.Range("C2").Value = 20

$priceSheet.Range("A3").Value = "2020-01-01T01:00:00"
$priceSheet.Range("B3").Value = "Base"
$priceSheet.Range("C3").Value = 20

$priceSheet.Range("A4").Value = "2020-01-01T02:00:00"
$priceSheet.Range("B4").Value = "Base"
$priceSheet.Range("C4").Value = 20

$priceSheet.Range("A5").Value = "2020-01-01T03:00:00"
$priceSheet.Range("B5").Value = "Base"
$priceSheet.Range("C5").Value = 20

$priceSheet.Range("A6").Value = "2020-01-01T04:00:00"
$priceSheet.Range("B6").Value = "Base"
$priceSheet.Range("C6").Value = 20
